$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Save" header in H1 and mirror the header style used by the
# other header cells (B1:G1) by copying formats from G1 (xlPasteFormats)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column values: 0 for rows 2-10, 1 for rows 11-13
$saveValues = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 1, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
